# Add completed triangle test plan (xlsx)
# Fills in the Method Inputs / Condition being Tested / Expected Result
# columns (E, F, G) for rows 7-14 of the test plan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$apostrophe = [char]0x2019

# Row 7 - __init__ happy path
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = 'color="red", side_1=5, side_2=6, side_3=7'
$ws.Range("G7").Value = "Triangle object created successfully, attributes set to input values"

# Method Inputs for the remaining __init__ failure-case rows are also "None"
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"

# __str__ / calculate_area / calculate_perimeter preconditions
$ws.Range("E12").Value = 'Valid triangle exists (e.g., color="red", sides=5,6,7)'
$ws.Range("E13").Value = 'Valid triangle exists (e.g., color="blue", sides=3,4,5)'
$ws.Range("E14").Value = 'Valid triangle exists (e.g., color="green", sides=2,3,4)'

# Remaining __init__ inputs (failure cases)
$ws.Range("F9").Value = 'color="red", side_1=5.2, side_2=6, side_3=7'
$ws.Range("F8").Value = 'color=" ", side_1=5, side_2=6, side_3=7'
$ws.Range("F10").Value = 'color="red", side_1=5, side_2="6", side_3=7'
$ws.Range("F11").Value = 'color="red", side_1=5, side_2=6, side_3=None'

# Actions for __str__ / calculate_area / calculate_perimeter
$ws.Range("F12").Value = "Call str(triangle)"
$ws.Range("F13").Value = "Call calculate_area()"
$ws.Range("F14").Value = "Call calculate_perimeter()"

# __init__ failure expected results
$ws.Range("G8").Value = 'ValueError("Color cannot be blank.")'
$ws.Range("G9").Value = 'ValueError("Side 1 must be numeric.")'
$ws.Range("G10").Value = 'ValueError("Side 2 must be numeric.")'
$ws.Range("G11").Value = 'ValueError("Side 3 must be numeric.")'

# Expected results for __str__ / calculate_area / calculate_perimeter
$ws.Range("G12").Value = 'String contains both "The shape color is red." and "5, 6 and 7"'
$ws.Range("G13").Value = "Returns 6.0 (area calculated by Heron${apostrophe}s formula)"
$ws.Range("G14").Value = "Returns 9.0"
